$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bell ring mode: add new question rows 27-54, mirroring the existing bank columns ---

# Row 27
$ws.Cells.Item(27, 1).Value = "{a} + {b} = "
$ws.Cells.Item(27, 2).Value = "Bellring"
$ws.Cells.Item(27, 3).Value = "a1:9*b1:9*"
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = "{a}+{b}"
$ws.Cells.Item(27, 6).Value = 10
$ws.Cells.Item(27, 7).Value = "रामू के पास {x} सेब हैं"

# Row 28
$ws.Cells.Item(28, 1).Value = "{a} - {b} = "
$ws.Cells.Item(28, 2).Value = "Bellring"
$ws.Cells.Item(28, 3).Value = "a5:9*b1:4*"
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = "{a}-{b}"
$ws.Cells.Item(28, 6).Value = 10
$ws.Cells.Item(28, 7).Value = "{a}-{b}"

# Row 29
$ws.Cells.Item(29, 1).Value = "{a} x {b} = "
$ws.Cells.Item(29, 2).Value = "Bellring"
$ws.Cells.Item(29, 3).Value = "a2:5*b2:5*"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).Value = "{a}*{b}"
$ws.Cells.Item(29, 6).Value = 10
$ws.Cells.Item(29, 7).Value = "{a}*{b}"

# Row 30
$ws.Cells.Item(30, 1).Value = "{a} / {b} = "
$ws.Cells.Item(30, 2).Value = "Bellring"
$ws.Cells.Item(30, 3).Value = "a2;1;4*b2;1;1*"
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = "{a}/{b}"
$ws.Cells.Item(30, 6).Value = 10
$ws.Cells.Item(30, 7).Value = "{a}/{b}"

# Row 31
$ws.Cells.Item(31, 1).Value = "{a} + {b} = "
$ws.Cells.Item(31, 2).Value = "Bellring"
$ws.Cells.Item(31, 3).Value = "a2;1;4*b1:5*"
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(31, 5).Value = "{a}+{b}"
$ws.Cells.Item(31, 6).Value = 10
$ws.Cells.Item(31, 7).Value = "{a}+{b}"

# Row 32
$ws.Cells.Item(32, 1).Value = "{a} x {b} = "
$ws.Cells.Item(32, 2).Value = "Bellring"
$ws.Cells.Item(32, 3).Value = "a3;1;3*b1:3*"
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = "{a}*{b}"
$ws.Cells.Item(32, 6).Value = 10
$ws.Cells.Item(32, 7).Value = "{a}*{b}"

# Row 33
$ws.Cells.Item(33, 1).Value = "{a} + {b} ="
$ws.Cells.Item(33, 2).Value = "Bellring"
$ws.Cells.Item(33, 3).Value = "a5;1;1*b1:9*"
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = "{a}+{b}"
$ws.Cells.Item(33, 6).Value = 10
$ws.Cells.Item(33, 7).Value = "{a}+{b}"

# Row 34
$ws.Cells.Item(34, 1).Value = "{a} + {b} ="
$ws.Cells.Item(34, 2).Value = "Bellring"
$ws.Cells.Item(34, 3).Value = "a1:5*b1:5*"
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(34, 5).Value = "{a}+{b}"
$ws.Cells.Item(34, 6).Value = 10
$ws.Cells.Item(34, 7).Value = "{a}+{b}"

# Row 35
$ws.Cells.Item(35, 1).Value = "{a} + {b} ="
$ws.Cells.Item(35, 2).Value = "Bellring"
$ws.Cells.Item(35, 3).Value = "a1:9*b1:1*"
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = "{a}+{b}"
$ws.Cells.Item(35, 6).Value = 10
$ws.Cells.Item(35, 7).Value = "{a}+{b}"

# Row 36
$ws.Cells.Item(36, 1).Value = "{a} + {b} ="
$ws.Cells.Item(36, 2).Value = "Bellring"
$ws.Cells.Item(36, 3).Value = "a1:7*b2:2*"
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 5).Value = "{a}+{b}"
$ws.Cells.Item(36, 6).Value = 10
$ws.Cells.Item(36, 7).Value = "{a}+{b}"

# Row 37
$ws.Cells.Item(37, 1).Value = "{a} + {b} ="
$ws.Cells.Item(37, 2).Value = "Bellring"
$ws.Cells.Item(37, 3).Value = "a2;1;4*b2;1;4*"
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 5).Value = "{a}+{b}"
$ws.Cells.Item(37, 6).Value = 10
$ws.Cells.Item(37, 7).Value = "{a}+{b}"

# Row 38
$ws.Cells.Item(38, 1).Value = "{a} + {b} ="
$ws.Cells.Item(38, 2).Value = "Bellring"
$ws.Cells.Item(38, 3).Value = "a3:3*b1:6*"
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(38, 5).Value = "{a}+{b}"
$ws.Cells.Item(38, 6).Value = 10
$ws.Cells.Item(38, 7).Value = "{a}+{b}"

# Row 39
$ws.Cells.Item(39, 1).Value = "{a} - {b} ="
$ws.Cells.Item(39, 2).Value = "Bellring"
$ws.Cells.Item(39, 3).Value = "a3:7*b1:2*"
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(39, 5).Value = "{a}-{b}"
$ws.Cells.Item(39, 6).Value = 10
$ws.Cells.Item(39, 7).Value = "{a}-{b}"

# Row 40
$ws.Cells.Item(40, 1).Value = "{a} - {b} ="
$ws.Cells.Item(40, 2).Value = "Bellring"
$ws.Cells.Item(40, 3).Value = "a2:9*b1:1*"
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(40, 5).Value = "{a}-{b}"
$ws.Cells.Item(40, 6).Value = 10
$ws.Cells.Item(40, 7).Value = "{a}-{b}"

# Row 41
$ws.Cells.Item(41, 1).Value = "{a} - {b} ="
$ws.Cells.Item(41, 2).Value = "Bellring"
$ws.Cells.Item(41, 3).Value = "a6:9*b1:5*"
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = "{a}-{b}"
$ws.Cells.Item(41, 6).Value = 10
$ws.Cells.Item(41, 7).Value = "{a}-{b}"

# Row 42
$ws.Cells.Item(42, 1).Value = "{a} - {b} ="
$ws.Cells.Item(42, 2).Value = "Bellring"
$ws.Cells.Item(42, 3).Value = "a3:9*b2:2*"
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(42, 5).Value = "{a}-{b}"
$ws.Cells.Item(42, 6).Value = 10
$ws.Cells.Item(42, 7).Value = "{a}-{b}"

# Row 43
$ws.Cells.Item(43, 1).Value = "{a} - {b} ="
$ws.Cells.Item(43, 2).Value = "Bellring"
$ws.Cells.Item(43, 3).Value = "a9:9*b1:8*"
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(43, 5).Value = "{a}-{b}"
$ws.Cells.Item(43, 6).Value = 10
$ws.Cells.Item(43, 7).Value = "{a}-{b}"

# Row 44
$ws.Cells.Item(44, 1).Value = "{a} x {b} ="
$ws.Cells.Item(44, 2).Value = "Bellring"
$ws.Cells.Item(44, 3).Value = "a2:2*b1:4*"
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(44, 5).Value = "{a}*{b}"
$ws.Cells.Item(44, 6).Value = 10
$ws.Cells.Item(44, 7).Value = "{a}*{b}"

# Row 45
$ws.Cells.Item(45, 1).Value = "{a} x {b} ="
$ws.Cells.Item(45, 2).Value = "Bellring"
$ws.Cells.Item(45, 3).Value = "a3:3*b1:3*"
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(45, 5).Value = "{a}*{b}"
$ws.Cells.Item(45, 6).Value = 10
$ws.Cells.Item(45, 7).Value = "{a}*{b}"

# Row 46
$ws.Cells.Item(46, 1).Value = "{a} x {b} ="
$ws.Cells.Item(46, 2).Value = "Bellring"
$ws.Cells.Item(46, 3).Value = "a1:3*b1:3*"
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = "{a}*{b}"
$ws.Cells.Item(46, 6).Value = 10
$ws.Cells.Item(46, 7).Value = "{a}*{b}"

# Row 47
$ws.Cells.Item(47, 1).Value = "{a} x {b} ="
$ws.Cells.Item(47, 2).Value = "Bellring"
$ws.Cells.Item(47, 3).Value = "a1:9*b1:1*"
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(47, 5).Value = "{a}*{b}"
$ws.Cells.Item(47, 6).Value = 10
$ws.Cells.Item(47, 7).Value = "{a}*{b}"

# Row 48
$ws.Cells.Item(48, 1).Value = "{a} x {b} ="
$ws.Cells.Item(48, 2).Value = "Bellring"
$ws.Cells.Item(48, 3).Value = "a4:4*b1:2*"
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(48, 5).Value = "{a}*{b}"
$ws.Cells.Item(48, 6).Value = 10
$ws.Cells.Item(48, 7).Value = "{a}*{b}"

# Row 49
$ws.Cells.Item(49, 1).Value = "{a} x {b} ="
$ws.Cells.Item(49, 2).Value = "Bellring"
$ws.Cells.Item(49, 3).Value = "a5:5*b1:1*"
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(49, 5).Value = "{a}*{b}"
$ws.Cells.Item(49, 6).Value = 10
$ws.Cells.Item(49, 7).Value = "{a}*{b}"

# Row 50
$ws.Cells.Item(50, 1).Value = "{a} / {b} ="
$ws.Cells.Item(50, 2).Value = "Bellring"
$ws.Cells.Item(50, 3).Value = "a2;1;4*b2:2*"
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(50, 5).Value = "{a}/{b}"
$ws.Cells.Item(50, 6).Value = 10
$ws.Cells.Item(50, 7).Value = "{a}/{b}"

# Row 51
$ws.Cells.Item(51, 1).Value = "{a} / {b} ="
$ws.Cells.Item(51, 2).Value = "Bellring"
$ws.Cells.Item(51, 3).Value = "a3;1;3*b3:3*"
$ws.Cells.Item(51, 4).Value = 1
$ws.Cells.Item(51, 5).Value = "{a}/{b}"
$ws.Cells.Item(51, 6).Value = 10
$ws.Cells.Item(51, 7).Value = "{a}/{b}"

# Row 52
$ws.Cells.Item(52, 1).Value = "{a} / {b} ="
$ws.Cells.Item(52, 2).Value = "Bellring"
$ws.Cells.Item(52, 3).Value = "a1:9*b1:1*"
$ws.Cells.Item(52, 4).Value = 1
$ws.Cells.Item(52, 5).Value = "{a}/{b}"
$ws.Cells.Item(52, 6).Value = 10
$ws.Cells.Item(52, 7).Value = "{a}/{b}"

# Row 53
$ws.Cells.Item(53, 1).Value = "{a} / {b} ="
$ws.Cells.Item(53, 2).Value = "Bellring"
$ws.Cells.Item(53, 3).Value = "a4;1;2*b4:4*"
$ws.Cells.Item(53, 4).Value = 1
$ws.Cells.Item(53, 5).Value = "{a}/{b}"
$ws.Cells.Item(53, 6).Value = 10
$ws.Cells.Item(53, 7).Value = "{a}/{b}"

# Row 54
$ws.Cells.Item(54, 1).Value = "{a} / {b} ="
$ws.Cells.Item(54, 2).Value = "Bellring"
$ws.Cells.Item(54, 3).Value = "a6:6*b1,2,3,6*"
$ws.Cells.Item(54, 4).Value = 1
$ws.Cells.Item(54, 5).Value = "{a}/{b}"
$ws.Cells.Item(54, 6).Value = 10
$ws.Cells.Item(54, 7).Value = "{a}/{b}"

# Apply the "Normal 2"-style formatting (Calibri 11, theme text color) used for the new rows' A:F columns
$ws.Range("A27:F54").Font.ThemeColor = 1

# Restore view state: zoom to 83%, scroll so row 23 is at top, select G32
$win = $excel.ActiveWindow
$win.Zoom = 83
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("G32").Select()

